$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 164.6
$ws.Range("I42").Value = 81.333336
$ws.Range("J42").Value = 200.28572
$ws.Range("K42").Value = 244.000008
$ws.Range("L42").Value = 600.85716
$ws.Range("M42").Value = -14.00000800000001
$ws.Range("N42").Value = -1060.85716

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14103.255
$ws.Range("I32").Value = 14364.583
$ws.Range("K32").Value = 14364.583
$ws.Range("M32").Value = -14077.583
$ws.Range("H61").Value = 2170.5278
$ws.Range("I61").Value = 1753.6923
$ws.Range("J61").Value = 2406.1304
$ws.Range("K61").Value = 1753.6923
$ws.Range("L61").Value = 2406.1304
$ws.Range("M61").Value = -1541.6923
$ws.Range("N61").Value = -2830.1304
$ws.Range("H97").Value = 836.1489
$ws.Range("I97").Value = 774.85364
$ws.Range("J97").Value = 1255
$ws.Range("K97").Value = 774.85364
$ws.Range("L97").Value = 1255
$ws.Range("M97").Value = -278.85364
$ws.Range("N97").Value = -2247
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H136").Value = 2170.5278
$ws.Range("I136").Value = 1753.6923
$ws.Range("J136").Value = 2406.1304
$ws.Range("K136").Value = 5261.0769
$ws.Range("L136").Value = 7218.3912
$ws.Range("M136").Value = -2711.0769
$ws.Range("N136").Value = -12318.3912

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2096.7144
$ws.Range("I94").Value = 2096.7144
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2096.7144
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1645.7144
$ws.Range("N94").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3879814
$ws.Range("I31").Value = 2109.1562
$ws.Range("J31").Value = 6177713
$ws.Range("K31").Value = 2109.1562
$ws.Range("L31").Value = 6177713
$ws.Range("M31").Value = -1814.1562
$ws.Range("N31").Value = -6178303
$ws.Range("H34").Value = 3879814
$ws.Range("I34").Value = 2109.1562
$ws.Range("J34").Value = 6177713
$ws.Range("K34").Value = 2109.1562
$ws.Range("L34").Value = 6177713
$ws.Range("M34").Value = -1907.1562
$ws.Range("N34").Value = -6178117
$ws.Range("H42").Value = 7750
$ws.Range("J42").Value = 7750
$ws.Range("L42").Value = 7750
$ws.Range("N42").Value = -8936
$ws.Range("H99").Value = 1953.8
$ws.Range("I99").Value = 1946.909
$ws.Range("J99").Value = 1962.2222
$ws.Range("K99").Value = 1946.909
$ws.Range("L99").Value = 1962.2222
$ws.Range("M99").Value = -448.9090000000001
$ws.Range("N99").Value = -4958.2222
$ws.Range("H126").Value = 1953.8
$ws.Range("I126").Value = 1946.909
$ws.Range("J126").Value = 1962.2222
$ws.Range("K126").Value = 5840.727000000001
$ws.Range("L126").Value = 5886.6666
$ws.Range("M126").Value = -3370.727000000001
$ws.Range("N126").Value = -10826.6666

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 500126
$ws.Range("I20").Value = 250
$ws.Range("K20").Value = 750
$ws.Range("M20").Value = -523
$ws.Range("H68").Value = 1443.7821
$ws.Range("I68").Value = 1185.9445
$ws.Range("J68").Value = 1521.1333
$ws.Range("K68").Value = 3557.8335
$ws.Range("L68").Value = 4563.3999
$ws.Range("M68").Value = -2746.8335
$ws.Range("N68").Value = -6185.3999
$ws.Range("H71").Value = 1443.7821
$ws.Range("I71").Value = 1185.9445
$ws.Range("J71").Value = 1521.1333
$ws.Range("K71").Value = 10673.5005
$ws.Range("L71").Value = 13690.1997
$ws.Range("M71").Value = -6617.5005
$ws.Range("N71").Value = -21802.1997
$ws.Range("H107").Value = 9011.799999999999
$ws.Range("I107").Value = 8054.5386
$ws.Range("J107").Value = 10048.833
$ws.Range("K107").Value = 24163.6158
$ws.Range("L107").Value = 30146.499
$ws.Range("M107").Value = -22243.6158
$ws.Range("N107").Value = -33986.499
$ws.Range("H121").Value = 310444.9
$ws.Range("I121").Value = 365
$ws.Range("J121").Value = 517164.84
$ws.Range("K121").Value = 1095
$ws.Range("L121").Value = 1551494.52
$ws.Range("M121").Value = 215
$ws.Range("N121").Value = -1554114.52
$ws.Range("H124").Value = 1278.2
$ws.Range("J124").Value = 1278.2
$ws.Range("L124").Value = 3834.6
$ws.Range("N124").Value = -13654.6
$ws.Range("H131").Value = 4494.0625
$ws.Range("J131").Value = 1653.6522
$ws.Range("L131").Value = 4960.9566
$ws.Range("N131").Value = -15040.9566
$ws.Range("H137").Value = 12195.096
$ws.Range("I137").Value = 3193.6365
$ws.Range("J137").Value = 22096.7
$ws.Range("K137").Value = 9580.9095
$ws.Range("L137").Value = 66290.10000000001
$ws.Range("M137").Value = -4480.9095
$ws.Range("N137").Value = -76490.10000000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3540.1853
$ws.Range("I68").Value = 3417.5881
$ws.Range("J68").Value = 3748.6
$ws.Range("K68").Value = 3417.5881
$ws.Range("L68").Value = 3748.6
$ws.Range("M68").Value = -2668.5881
$ws.Range("N68").Value = -5246.6
$ws.Range("H71").Value = 3540.1853
$ws.Range("I71").Value = 3417.5881
$ws.Range("J71").Value = 3748.6
$ws.Range("K71").Value = 17087.9405
$ws.Range("L71").Value = 18743
$ws.Range("M71").Value = -13343.9405
$ws.Range("N71").Value = -26231
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 5750
$ws.Range("J100").Value = 5750
$ws.Range("L100").Value = 5750
$ws.Range("N100").Value = -6832

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 425122.2
$ws.Range("I136").Value = 497130.47
$ws.Range("K136").Value = 1491391.41
$ws.Range("M136").Value = -1488841.41

